# Part 1: add the new "Today's Attendance password / compare" text box to slide 1
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.AddTextbox(1, 240.00007874015748, -1.6296850393700788, 464.84220472440944, 82.39685039370079)
$shp.Name = "TextBox 2"

$tf = $shp.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "Today's Attendance password`rcompare"

$p1 = $tr.Paragraphs(1, 1)
$p1.Font.Highlight.RGB = 65535

$p2 = $tr.Paragraphs(2, 1)
$p2.Font.Size = 44
$p2.Font.Highlight.RGB = 65535

$shp.Height = 82.39685039370079

$shp.Fill.Visible = $false
$shp.Line.Visible = $false
$shp.Line.Weight = 1
